$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values (rows 2-5)
$ws.Range("A2").Value = 12
$ws.Range("B2").Value = 8

$ws.Range("A3").Value = 11
$ws.Range("B3").Value = 3

$ws.Range("A4").Value = 22
$ws.Range("B4").Value = 1

$ws.Range("A5").Value = 21
$ws.Range("B5").Value = 1

# Remove row 6 (was A6=112, B6=1)
$ws.Range("A6:B6").Delete()
